
$d = $word.ActiveDocument

# --- 1. Title paragraph: re-split runs " para" | " el s" | "istema de " | Tickets | " soporte..."
#        into " para el s" | "istema de " | Tickets | " soporte..." (one fewer run),
#        moving the gramStart/gramEnd proofErr markers to still bracket "Tickets".
$titlePara = $d.Paragraphs.Item(3)
$titleRange = $titlePara.Range
$findTitle = $titleRange.Duplicate
$findTitle.Find.Execute(" para", $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
$titleTarget = $d.Range($findTitle.Start, $titlePara.Range.End - 1)
$titleXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> para el s</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve">istema de </w:t></w:r><w:proofErr w:type="gramStart"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="es-ES"/></w:rPr><w:t>Tickets</w:t></w:r><w:proofErr w:type="gramEnd"/><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="32"/><w:szCs w:val="32"/><w:lang w:val="es-ES"/></w:rPr><w:t xml:space="preserve"> soporte técnico para el Centro de Cómputo del ITL</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$titleTarget.InsertXML($titleXml)

# --- 2. Paragraph "La pagina web..." : change "el cuantas consultas" -> "las consultas"
#        and drop the gramStart/gramEnd proofErr markers around it.
$p14 = $d.Paragraphs.Item(14)
$p14Range = $p14.Range
$p14Target = $d.Range($p14Range.Start, $p14Range.End - 1)
$p14Xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve">La pagina web, base de datos y API, al poder estar alojadas en un servidor para su despliegue, dependerá de que tan potente sea el servidor en recibir a los clientes </w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t>las consultas</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="24"/><w:szCs w:val="24"/></w:rPr><w:t xml:space="preserve"> podrá manejar el formulario web, ya que las entradas de la API, al estar conectada a una base de datos de SQL SERVER aguantan demasiada demanda, por lo que, esta parte dependerá meramente de la calidad del servidor que tenga el Instituto Tecnológico de León.</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p14Target.InsertXML($p14Xml)

# --- 3. Wrap the "El sistema..." through "La pagina web..." paragraphs in a bookmark.
$p8 = $d.Paragraphs.Item(8)
$p14b = $d.Paragraphs.Item(14)
$bkRange = $d.Range($p8.Range.Start, $p14b.Range.End - 1)
$d.Bookmarks.Add("_Hlk167752534", $bkRange)

Write-Host "All edits applied"
